# Final cleanup pass on the test-account-password sheet:
#  - drop the scratch "JSON blob" helper column (C) along with its formulas
#  - re-sync rows 17-18 with the plain data-row formatting (their old
#    fill/border styling existed only to go with the now-removed column C)
#  - upload the latest export from the database: a few newly-created
#    sample accounts, plus some blank (but formatted) trailing rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Remove the helper "JSON" column C's contents (values + formulas),
#    leaving the column's width formatting alone.
#    This also drops the now-unused "Json" shared string automatically.
# ---------------------------------------------------------------------
$ws.Range("C1:C18").Clear()

# ---------------------------------------------------------------------
# 2) Append the newly-created sample accounts from the latest DB export.
# ---------------------------------------------------------------------
$ws.Range("A19").Value = "jhess1"
$ws.Range("B19").Value = "examplePatient"

$ws.Range("A20").Value = "hphilips1"
$ws.Range("B20").Value = "exampleDoctor"

$ws.Range("A21").Value = "ahill1"
$ws.Range("B21").Value = "exampleAdmin"

# ---------------------------------------------------------------------
# 3) Rows 17-18 (qcee1 / ecee1) keep their username+password values, but
#    should now look like the rest of the plain data rows instead of the
#    special highlighted style they had next to the JSON column. Rows
#    19-21 (the new accounts above) get the same plain data-row look.
# ---------------------------------------------------------------------
$ws.Range("A2").Copy()
$ws.Range("A17:A21").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("B2").Copy()
$ws.Range("B17:B21").PasteSpecial(-4122)   # xlPasteFormats

# ---------------------------------------------------------------------
# 4) Pad a few trailing formatted-but-empty rows, matching the rest of
#    the freshly pasted range from the database export.
# ---------------------------------------------------------------------
$ws.Range("A2").Copy()
$ws.Range("A22:A27").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("B2").Copy()
$ws.Range("B22:B27").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 5) Leave the selection where the cursor ended up after the edit.
# ---------------------------------------------------------------------
$ws.Range("C14").Select()
